# "measure diameter for normal temperature"
# Fill in the newly-measured diameter column (C) on the "常温" (normal
# temperature) sheet, then leave the workbook focused/scrolled the way the
# author left it (常温 sheet active instead of 常温硬度).

$wb = $excel.ActiveWorkbook

# Sheets, in tab order: 1=冷库, 2=常温, 3=冷库硬度, 4=常温硬度
$wsNormal = $wb.Worksheets.Item(2)

# New column-C measurements for rows 37 and 41-99 on the 常温 sheet.
# (Rows 38-40 were left unmeasured - only their row span metadata changes.)
$cValues = @{
    37 = 46.75
    41 = 47.91
    42 = 43.89
    43 = 43.86
    44 = 47.91
    45 = 43.23
    46 = 45.97
    47 = 49.55
    48 = 49.95
    49 = 49.74
    50 = 42.17
    51 = 47.97
    52 = 49.67
    53 = 45.96
    54 = 47.78
    55 = 47.17
    56 = 45.53
    57 = 47.23
    58 = 46.65
    59 = 42.71
    60 = 46.56
    61 = 48.09
    62 = 47.6
    63 = 48.41
    64 = 46.63
    65 = 42.06
    66 = 45.49
    67 = 46.29
    68 = 49.71
    69 = 47.05
    70 = 44.25
    71 = 45.12
    72 = 48.34
    73 = 47.09
    74 = 50.05
    75 = 49.44
    76 = 48.61
    77 = 47.83
    78 = 50.06
    79 = 46.04
    80 = 45.77
    81 = 49.06
    82 = 47.4
    83 = 47.17
    84 = 48.06
    85 = 48.48
    86 = 48.11
    87 = 49.5
    88 = 48.83
    89 = 49.26
    90 = 48.42
    91 = 45.38
    92 = 49.56
    93 = 44.59
    94 = 47.59
    95 = 45.41
    96 = 50.35
    97 = 47.26
    98 = 47.23
    99 = 51.67
}

foreach ($row in $cValues.Keys) {
    $wsNormal.Cells.Item($row, 3).Value = $cValues[$row]
}

# The 常温硬度 sheet (previously the active tab) also scrolled slightly
# while the author was working, before focus moved away from it.
$wsHardness = $wb.Worksheets.Item(4)
$wsHardness.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1

# Switch focus to the 常温 sheet (it becomes the active/selected tab,
# replacing 常温硬度), scroll it and select C67, matching the author's
# final view state.
$wsNormal.Activate()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$wsNormal.Range("C67").Select()
